# Update technical factor return statistics on Sheet1 with refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.362096774193548, 1.530645161290322, -2.050806451612903, 9.741935483870968, -2.245161290322581, -3.092741935483871)
    3  = @(0.4548387096774194, 2.263709677419354, 4.882258064516127, 19.81209677419354, 11.59838709677419, 3.642741935483873)
    4  = @(1.203225806451613, 1.265322580645161, -1.523387096774194, 11.20564516129032, -1.191129032258064, -2.838709677419355)
    5  = @(0.428225806451613, 2.819354838709677, 4.926612903225809, 18.78629032258064, 13.07661290322581, 4.187096774193548)
    6  = @(0.1338709677419355, 1.532258064516129, 4.509677419354838, 14.94838709677419, 13.84354838709677, 4.733064516129032)
    7  = @(1.700806451612903, 3.114516129032256, 0.8137096774193551, 15.73467741935484, -0.5782258064516129, -1.450806451612904)
    8  = @(1.3688, 3.395999999999999, 4.1232, 19.8624, 7.4296, 1.8144)
    9  = @(0.2792, 1.332800000000001, 1.617599999999999, 10.112, 7.529599999999998, 1.7832)
    10 = @(0.7536, 2.4064, 3.3104, 24.3456, 18.228, 2.768800000000001)
    11 = @(1.0912, 2.159200000000001, 0.7479999999999987, 3.375999999999997, -7.343548387096773, -1.2976)
    12 = @(1.0616, 3.4856, 5.825600000000001, 24.316, 13.1248, 2.9248)
    13 = @(0.62, 0.7199999999999999, -0.5568, 5.396, 0.7701612903225794, 0.3288000000000005)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i  # column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
